# Updates cryptos list - price/volume refresh + a handful of row reshuffles.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "66.811.66"
$ws.Range("E2").Value = "  +0.47%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "2.518.21"
$ws.Range("E3").Value = "  -2.39%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.03%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "587.32"
$ws.Range("E5").Value = "  +0.82%  "

# --- Row 6 (Solana) ---
$ws.Range("D6").Value = "171.38"
$ws.Range("E6").Value = "  +3.34%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.05%  "

# --- Row 8 (XRP) ---
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.47%  "

# --- Row 9 (LidoStakedEther) ---
$ws.Range("D9").Value = "2.517.84"
$ws.Range("E9").Value = "  -2.35%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  +0.03%  "

# --- Row 11 (TRON) ---
$ws.Range("E11").Value = "  +1.68%  "

# --- Row 12 (Toncoin) ---
$ws.Range("E12").Value = "  -1.27%  "

# --- Row 13 (Cardano) ---
$ws.Range("D13").Value = "0.341"

# --- Row 14 (Avalanche) ---
$ws.Range("D14").Value = "26.45"
$ws.Range("E14").Value = "  -1.33%  "

# --- Row 15 (WrappedliquidstakedEther2.0) ---
$ws.Range("D15").Value = "2.982.06"
$ws.Range("E15").Value = "  -2.26%  "

# --- Row 16 (ShibaInu) ---
$ws.Range("E16").Value = "  -1.59%  "

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "66.687.82"
$ws.Range("E17").Value = "  +0.52%  "

# --- Row 18 (WrappedEther) ---
$ws.Range("D18").Value = "2.516.52"
$ws.Range("E18").Value = "  -2.46%  "

# --- Row 19 (Uniswap) ---
$ws.Range("E19").Value = "  +2.63%  "

# --- Row 20 (Chainlink) ---
$ws.Range("D20").Value = "11.26"
$ws.Range("E20").Value = "  -1.50%  "

# --- Row 21 (BitcoinCash) ---
$ws.Range("E21").Value = "  +0.27%  "

# --- Row 22 (Polkadot) ---
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -1.74%  "

# --- Row 23 (NEARProtocol) ---
$ws.Range("E23").Value = "  -0.76%  "

# --- Row 24 (SuiNetwork) ---
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +4.70%  "

# --- Row 25 (Dai) ---
$ws.Range("E25").Value = "  +0.02%  "

# --- Row 26 (Litecoin) ---
$ws.Range("D26").Value = "69.68"
$ws.Range("E26").Value = "  +1.01%  "

# --- Row 27 (Aptos) ---
$ws.Range("E27").Value = "  -0.50%  "

# --- Row 28 / 29 swap (WrappedeETH <-> Binance-PegBSC-USD) with new values ---
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.659.28"
$ws.Range("E29").Value = "  -1.89%  "

# --- Row 30 (PEPE) ---
$ws.Range("E30").Value = "  -1.75%  "

# --- Row 31 (Bittensor) ---
$ws.Range("D31").Value = "530.56"
$ws.Range("E31").Value = "  -1.20%  "

# --- Row 32 (InternetComputer(DFINITY)) ---
$ws.Range("E32").Value = "  +0.55%  "

# --- Row 33 (Fetch.AI) ---
$ws.Range("E33").Value = "  -0.31%  "

# --- Row 34 (PancakeSwap) ---
$ws.Range("E34").Value = "  -0.96%  "

# --- Row 35 (Kaspa) ---
$ws.Range("E35").Value = "  -1.23%  "

# --- Row 36 (FirstDigitalUSD) ---
$ws.Range("E36").Value = "  +0.09%  "

# --- Row 37 / 38 swap (ImmutableX <-> Monero) with new values ---
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "157.86"
$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -1.14%  "

# --- Row 39 (EthereumClassic) ---
$ws.Range("E39").Value = "  -1.25%  "

# --- Row 40 (WhiteBITCoin) ---
$ws.Range("D40").Value = "18.42"
$ws.Range("E40").Value = "  +1.03%  "

# --- Row 42 (Stacks) ---
$ws.Range("E42").Value = "  -0.89%  "

# --- Row 43 (RenderToken) ---
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  -0.95%  "

# --- Row 44 (USDe) ---
$ws.Range("E44").Value = "  -0.02%  "

# --- Row 45 (dogwifhat) ---
$ws.Range("E45").Value = "  +2.26%  "

# --- Row 46 / 47 / 48 / 49 / 50 / 51 cascade shift (OKB, Aave, ARBITRUM, BabyDogeCoin, Filecoin, Optimism -> Aave, ARBITRUM, BabyDogeCoin, Filecoin, Optimism, Cronos) ---
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "148.47"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.554"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0275"
$ws.Range("E48").Value = "  -4.65%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "3.66"
$ws.Range("E49").Value = "  -1.76%  "

$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0755"
$ws.Range("E51").Value = "  -0.50%  "
